$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles, row height) of row 3 into the new row 5 so
# that the new row matches the look of the other "Observations" rows.
$ws.Range("A3:D3").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)
$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(3).RowHeight

# Populate the new row with the rolled-back CodeBaseRegister note.
$ws.Range("A5").Value = 43401.843055555553
$ws.Range("B5").Value = 0.0098379629629629633
$ws.Range("C5").Value = "643,240.11 KB"
$ws.Range("D5").Value = "Evaluated one individual on desktop with debug code single thread (break point at PushPG.compute_errors() line #38)."

# Match the author's final selection in the saved workbook.
$ws.Range("C5").Select()
